# "set up test for table filter, to use between adding and installing games"
#
# The WBEpicSheet / WBSteamSheet config rows (and their trailing blank
# separator row) are removed from the Settings sheet, and the description
# for the EpicCredential setting is reworded to call out that it is a
# *local* credential (to distinguish it from the table-filter test step
# that sits between the "adding" and "installing" game steps).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Remove the WBEpicSheet / WBSteamSheet rows (A32:A34 -> "WBEpicSheet",
# "WBSteamSheet", <blank>), shifting everything below up by 3 rows.
$ws.Range("A32:A34").EntireRow.Delete()

# The EpicCredential row (previously row 35, now row 32) gets a reworded
# description.
$ws.Range("C32").Value = "Name for local credential to sign into Epic Games"

# Leave the selection where the edit happened, scrolled back to the top.
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("C32").Select()
